$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.316.48'
$ws.Range("E2").Value = '  -1.31%  '
$ws.Range("D3").Value = '3.553.44'
$ws.Range("E3").Value = '  +1.08%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.24'
$ws.Range("E5").Value = '  -0.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.47'
$ws.Range("E6").Value = '  -2.54%  '
$ws.Range("D7").Value = '3.552.09'
$ws.Range("E7").Value = '  +1.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.481'
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("B10").Value = 'Toncoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '8.08'
$ws.Range("E10").Value = '  +0.67%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.137'
$ws.Range("E11").Value = '  -4.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.411'
$ws.Range("E12").Value = '  -2.73%  '
$ws.Range("D13").Value = '4.155.26'
$ws.Range("E13").Value = '  +1.13%  '
$ws.Range("E14").Value = '  -4.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '30.12'
$ws.Range("E15").Value = '  -4.45%  '
$ws.Range("D16").Value = '3.556.19'
$ws.Range("E16").Value = '  +1.29%  '
$ws.Range("D17").Value = '66.372.17'
$ws.Range("E17").Value = '  -1.25%  '
$ws.Range("E18").Value = '  -0.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.13'
$ws.Range("E19").Value = '  +1.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.21'
$ws.Range("E20").Value = '  -2.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.94'
$ws.Range("E21").Value = '  -3.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '427.98'
$ws.Range("E22").Value = '  -2.06%  '
$ws.Range("E23").Value = '  -1.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.93'
$ws.Range("E24").Value = '  -1.32%  '
$ws.Range("D25").Value = '3.692.10'
$ws.Range("E25").Value = '  +1.02%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").Value = '  +1.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.12'
$ws.Range("E28").Value = '  -1.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.20'
$ws.Range("E29").Value = '  -6.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.49'
$ws.Range("E30").Value = '  -1.39%  '
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("E32").Value = '  -5.62%  '
$ws.Range("E33").Value = '  -4.27%  '
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.41'
$ws.Range("E34").Value = '  -0.82%  '
$ws.Range("B35").Value = 'RenzoRestakedETH'
$ws.Range("C35").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D35").Value = '3.544.40'
$ws.Range("E35").Value = '  +1.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.76'
$ws.Range("E37").Value = '  -2.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.65'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.83'
$ws.Range("E39").Value = '  -2.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '174.65'
$ws.Range("E41").Value = '  -0.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0859'
$ws.Range("E42").Value = '  -4.52%  '
$ws.Range("E43").Value = '  -2.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.895'
$ws.Range("E44").Value = '  -0.30%  '
$ws.Range("E45").Value = '  -6.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '45.77'
$ws.Range("E46").Value = '  -1.28%  '
$ws.Range("E47").Value = '  -1.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '26.03'
$ws.Range("E48").Value = '  -7.57%  '
$ws.Range("E49").Value = '  -2.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.14'
$ws.Range("E50").Value = '  -4.49%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.91'
$ws.Range("E51").Value = '  +5.31%  '
